$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 394. This shifts the former rows
# 394..482 down to 396..484 (they keep all their original values/styles),
# matching the diff's row-by-row "D/J/K/L/M/O/P shifted down by two" pattern,
# and matching the new dimension A1:R484 plus the two brand-new <row> blocks
# that appear right after the old row 482 in the diff (those correspond to
# what is now rows 483/484, i.e. the content of the former rows 481/482).
$ws.Rows("394:395").Insert()

# Populate the first new row (394) with its data.
$ws.Range("A394").Value = 10
$ws.Range("B394").Value = "Vega Modelo de Temuco"
$ws.Range("C394").Value = "La Araucanía"
$ws.Range("D394").Value = 44855
$ws.Range("E394").Value = 9
$ws.Range("F394").Value = 100112008
$ws.Range("G394").Value = "Coliflor"
$ws.Range("H394").Value = "Sin especificar"
$ws.Range("I394").Value = "Primera"
$ws.Range("J394").Value = 1200
$ws.Range("K394").Value = 1000
$ws.Range("L394").Value = 1000
$ws.Range("M394").Value = 1000
$ws.Range("N394").Value = "$/unidad"
$ws.Range("O394").Value = "Región Metropolitana"
$ws.Range("P394").Value = 1000
$ws.Range("Q394").Value = 1
$ws.Range("R394").Value = "Hortaliza"

# Populate the second new row (395) with its data.
$ws.Range("A395").Value = 10
$ws.Range("B395").Value = "Vega Modelo de Temuco"
$ws.Range("C395").Value = "La Araucanía"
$ws.Range("D395").Value = 44855
$ws.Range("E395").Value = 9
$ws.Range("F395").Value = 100112008
$ws.Range("G395").Value = "Coliflor"
$ws.Range("H395").Value = "Sin especificar"
$ws.Range("I395").Value = "Primera"
$ws.Range("J395").Value = 600
$ws.Range("K395").Value = 1200
$ws.Range("L395").Value = 1200
$ws.Range("M395").Value = 1200
$ws.Range("N395").Value = "$/unidad"
$ws.Range("O395").Value = "Región del Maule"
$ws.Range("P395").Value = 1200
$ws.Range("Q395").Value = 1
$ws.Range("R395").Value = "Hortaliza"
